$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "G2"
$ws.Range("B3").Value = "Test1"
$ws.Range("C3").Value = "Daily"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 45860
$ws.Range("F3").Value = 30

# Match the date-formatted style used by E2 on the new E3 cell
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat
